# Fruta / hortaliza, semanal
# Insert a new weekly record at row 530 (pushing the existing rows 530:568
# down to 531:569) and populate it with the new Mandarina / Murcott price
# observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 530:568 down one row, creating a blank row 530.
$ws.Rows(530).Insert()

# Populate the new row with the latest market observation.
$ws.Range("A530").Value = 5
$ws.Range("B530").Value = "Macroferia Regional de Talca"
$ws.Range("C530").Value = "Maule"
$ws.Range("D530").Value = 45223
$ws.Range("E530").Value = 7
$ws.Range("F530").Value = "Fruta"
$ws.Range("G530").Value = 100102
$ws.Range("H530").Value = "Cítricos"
$ws.Range("I530").Value = 100102004
$ws.Range("J530").Value = "Mandarina"
$ws.Range("K530").Value = "Murcott"
$ws.Range("L530").Value = "Primera"
$ws.Range("M530").Value = 340
$ws.Range("N530").Value = 8000
$ws.Range("O530").Value = 8000
$ws.Range("P530").Value = 8000
$ws.Range("Q530").Value = "`$/bandeja 18 kilos"
$ws.Range("R530").Value = "Región de O'Higgins"
$ws.Range("S530").Value = 444
$ws.Range("T530").Value = 18
